$d = $word.ActiveDocument

# --- Paragraph 1 formatting changes ---
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat

# Add a paragraph border (top/left/bottom/right) each with a 5-twip space,
# but no line style/weight/color specified (matches target OOXML: <w:top w:space="5"/> etc.)
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt)
$pf.LeftIndent = 11.25

# --- Paragraph 1 text changes ---
# Replace the placeholder ID text (first run) and drop the trailing
# space-only run that follows it.
$r1 = $d.Range(0, 36)
$r1.Text = "**ID__AFFARS_AFMC_PGI_5315_3C__ID**"

# After the text substitution the paragraph's text is 35 characters long,
# followed by the leftover space run at [35,36) and the paragraph mark.
$trailingSpace = $d.Range(35, 36)
$trailingSpace.Delete()
